# Insert a new data row at row 859 (pushing the existing rows 859..941 down
# to 860..942) and populate it with the new record. Excel's native row
# insert shifts the existing cell values/styles down automatically, which
# matches the target diff (dimension grows from A1:R941 to A1:R942).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(859).Insert()

$ws.Range("A859").Value = 11
$ws.Range("B859").Value = "Vega Monumental Concepción"
$ws.Range("C859").Value = "Bíobío"
$ws.Range("D859").Value = 45212
$ws.Range("E859").Value = 8
$ws.Range("F859").Value = 100112004
$ws.Range("G859").Value = "Cebolla"
$ws.Range("H859").Value = "Sin especificar"
$ws.Range("I859").Value = "2a (guarda)"
$ws.Range("J859").Value = 150
$ws.Range("K859").Value = 12000
$ws.Range("L859").Value = 12000
$ws.Range("M859").Value = 12000
$ws.Range("N859").Value = "`$/malla 16 kilos"
$ws.Range("O859").Value = "Región de O'Higgins"
$ws.Range("P859").Value = 750
$ws.Range("Q859").Value = 16
$ws.Range("R859").Value = "Hortaliza"
